# Capacita - aprovacao de cursos
# Replace the placeholder "Teste1/2/3" request rows with the real
# course-approval requests and extend the table from 3 data rows (2-4)
# to 6 data rows (2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column layout (row 1 headers, unchanged):
#   A Descricao | B Plano | C Iniciativa | D Prioridade
#   E Quantidade de Servidores | F Area Conhecimento | G Nivel
#   H Hora de Duração | I Turno | J Mês

# Extend formatting (style + row height) for the new rows 5-7 by copying
# row 2's formatting down, then overwrite every cell's value below.
$ws.Range("A2:J2").Copy($ws.Range("A5:J7"))
$ws.Rows.Item(5).RowHeight = 20
$ws.Rows.Item(6).RowHeight = 20
$ws.Rows.Item(7).RowHeight = 20

$rows = @(
    @{ A = "";                                    D = "Baixa"; E = 21; F = "Teoria da Comunicação";   G = "facil"; H = 1; I = "Matutino";   J = "Dezembro" },
    @{ A = "";                                    D = "Baixa"; E = 22; F = "Teoria da Comunicação";   G = "medio"; H = 5; I = "Vespertino"; J = "Junho" },
    @{ A = "bligs blaps";                         D = "Alta";  E = 22; F = "Sociologia Jurídica";     G = "facil"; H = 5; I = "Matutino";   J = "Fevereiro" },
    @{ A = "joop joops";                          D = "Baixa"; E = 22; F = "Teoria do Estado";        G = "facil"; H = 2; I = "Matutino";   J = "Dezembro" },
    @{ A = "";                                    D = "Baixa"; E = 22; F = "Teoria da Comunicação";   G = "medio"; H = 2; I = "Matutino";   J = "Março" },
    @{ A = "topicos em negocios internacionais "; D = "Média"; E = 22; F = "Negócios Internacionais"; G = "facil"; H = 1; I = "Matutino";   J = "Novembro" }
)

$r = 2
foreach ($row in $rows) {
    if ($row.A -eq "") {
        $ws.Cells.Item($r, 1).ClearContents()
    } else {
        $ws.Cells.Item($r, 1).Value = $row.A
    }
    $ws.Cells.Item($r, 2).Value = "bla"
    $ws.Cells.Item($r, 3).Value = "Senado"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $r++
}
